$wb = $excel.ActiveWorkbook

# Sheet ALC, row 15 (hunk idx 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 153.23
$ws.Range("I15").Value = 153.23
$ws.Range("K15").Value = 459.6899999999999
$ws.Range("M15").Value = -290.6899999999999

# Sheet ALC, row 121 (hunk idx 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2625.1482
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2625.1482
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 7875.444600000001
$ws.Range("N121").Value = -11369.4446

# Sheet ALC, row 125 (hunk idx 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2099.2104
$ws.Range("I125").Value = 1283
$ws.Range("J125").Value = 2316.8667
$ws.Range("K125").Value = 11547
$ws.Range("L125").Value = 20851.8003
$ws.Range("M125").Value = -9087
$ws.Range("N125").Value = -25771.8003

# Sheet ALC, row 129 (hunk idx 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 765.0625
$ws.Range("I129").Value = 541.6667
$ws.Range("J129").Value = 899.1
$ws.Range("K129").Value = 1625.0001
$ws.Range("L129").Value = 2697.3
$ws.Range("M129").Value = 3374.9999
$ws.Range("N129").Value = -12697.3

# Sheet ARM, row 32 (hunk idx 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7839.46
$ws.Range("I32").Value = 6585.6064
$ws.Range("J32").Value = 27483.166
$ws.Range("K32").Value = 6585.6064
$ws.Range("L32").Value = 27483.166
$ws.Range("M32").Value = -6298.6064
$ws.Range("N32").Value = -28057.166

# Sheet ARM, row 117 (hunk idx 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 20066.666
$ws.Range("J117").Value = 20066.666
$ws.Range("L117").Value = 20066.666
$ws.Range("N117").Value = -29244.666

# Sheet ARM, row 135 (hunk idx 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 49402.637
$ws.Range("J135").Value = 49402.637
$ws.Range("L135").Value = 49402.637
$ws.Range("N135").Value = -59542.637

# Sheet BSM, row 22 (hunk idx 7)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 235.65218
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 610
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 610
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -956

# Sheet BSM, row 64 (hunk idx 8)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2026.3636
$ws.Range("I64").Value = 5069
$ws.Range("J64").Value = 287.7143
$ws.Range("K64").Value = 5069
$ws.Range("L64").Value = 287.7143
$ws.Range("M64").Value = -4844
$ws.Range("N64").Value = -737.7143

# Sheet BSM, row 67 (hunk idx 9)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 2026.3636
$ws.Range("I67").Value = 5069
$ws.Range("J67").Value = 287.7143
$ws.Range("K67").Value = 5069
$ws.Range("L67").Value = 287.7143
$ws.Range("M67").Value = -4289
$ws.Range("N67").Value = -1847.7143

# Sheet BSM, row 105 (hunk idx 10)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2503.3333
$ws.Range("J105").Value = 2750
$ws.Range("L105").Value = 2750
$ws.Range("N105").Value = -6244

# Sheet CRP, row 16 (hunk idx 11)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 908.2105
$ws.Range("I16").Value = 916
$ws.Range("J16").Value = 866.6667
$ws.Range("K16").Value = 916
$ws.Range("L16").Value = 866.6667
$ws.Range("M16").Value = -629
$ws.Range("N16").Value = -1440.6667

# Sheet CRP, row 31 (hunk idx 12)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5211719
$ws.Range("I31").Value = 2054.2856
$ws.Range("J31").Value = 6670425
$ws.Range("K31").Value = 2054.2856
$ws.Range("L31").Value = 6670425
$ws.Range("M31").Value = -1759.2856
$ws.Range("N31").Value = -6671015

# Sheet CRP, row 34 (hunk idx 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5211719
$ws.Range("I34").Value = 2054.2856
$ws.Range("J34").Value = 6670425
$ws.Range("K34").Value = 2054.2856
$ws.Range("L34").Value = 6670425
$ws.Range("M34").Value = -1852.2856
$ws.Range("N34").Value = -6670829

# Sheet CRP, row 50 (hunk idx 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10084.5
$ws.Range("J50").Value = 10084.5
$ws.Range("L50").Value = 10084.5
$ws.Range("N50").Value = -11334.5

# Sheet CRP, row 51 (hunk idx 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8835.182000000001
$ws.Range("J51").Value = 9666.333000000001
$ws.Range("L51").Value = 9666.333000000001
$ws.Range("N51").Value = -11138.333

# Sheet CRP, row 58 (hunk idx 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2939.043
$ws.Range("I58").Value = 1341.0312
$ws.Range("J58").Value = 4284.737
$ws.Range("K58").Value = 1341.0312
$ws.Range("L58").Value = 4284.737
$ws.Range("M58").Value = -1138.0312
$ws.Range("N58").Value = -4690.737

# Sheet CRP, row 59 (hunk idx 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 14455.3
$ws.Range("J59").Value = 15444.125
$ws.Range("L59").Value = 15444.125
$ws.Range("N59").Value = -17734.125

# Sheet CRP, row 60 (hunk idx 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9337.454
$ws.Range("J60").Value = 9971.200000000001
$ws.Range("L60").Value = 9971.200000000001
$ws.Range("N60").Value = -10993.2

# Sheet CRP, row 61 (hunk idx 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 8835.182000000001
$ws.Range("J61").Value = 9666.333000000001
$ws.Range("L61").Value = 9666.333000000001
$ws.Range("N61").Value = -10362.333

# Sheet CRP, row 68 (hunk idx 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17710.334
$ws.Range("J68").Value = 18023.273
$ws.Range("L68").Value = 18023.273
$ws.Range("N68").Value = -19521.273

# Sheet ARM, row 70 (hunk idx 21)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28630

# Sheet CRP, row 71 (hunk idx 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17710.334
$ws.Range("J71").Value = 18023.273
$ws.Range("L71").Value = 54069.819
$ws.Range("N71").Value = -61557.819

# Sheet ARM, row 73 (hunk idx 23)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -30184

# Sheet CRP, row 74 (hunk idx 24)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 14244
$ws.Range("J74").Value = 16366.556
$ws.Range("L74").Value = 16366.556
$ws.Range("N74").Value = -18114.556

# Sheet CRP, row 77 (hunk idx 25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 14244
$ws.Range("J77").Value = 16366.556
$ws.Range("L77").Value = 49099.66800000001
$ws.Range("N77").Value = -57835.66800000001

# Sheet CRP, row 113 (hunk idx 26)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 908.2105
$ws.Range("I113").Value = 916
$ws.Range("J113").Value = 866.6667
$ws.Range("K113").Value = 916
$ws.Range("L113").Value = 866.6667
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -5206.6667

# Sheet CRP, row 122 (hunk idx 27)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3476412.8
$ws.Range("I122").Value = 4469163
$ws.Range("J122").Value = 1787.5
$ws.Range("K122").Value = 13407489
$ws.Range("L122").Value = 5362.5
$ws.Range("M122").Value = -13405039
$ws.Range("N122").Value = -10262.5

# Sheet CRP, row 136 (hunk idx 28)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2939.043
$ws.Range("I136").Value = 1341.0312
$ws.Range("J136").Value = 4284.737
$ws.Range("K136").Value = 4023.0936
$ws.Range("L136").Value = 12854.211
$ws.Range("M136").Value = -1473.0936
$ws.Range("N136").Value = -17954.211

# Sheet CUL, row 33 (hunk idx 29)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 31250262
$ws.Range("I33").Value = 45454830
$ws.Range("J33").Value = 210.8
$ws.Range("K33").Value = 272728980
$ws.Range("L33").Value = 1264.8
$ws.Range("M33").Value = -272728697
$ws.Range("N33").Value = -1830.8

# Sheet ALC, row 95 (hunk idx 30)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 2500
$ws.Range("J95").Value = 2500
$ws.Range("L95").Value = 7500
$ws.Range("N95").Value = -11618

# Sheet CUL, row 131 (hunk idx 31)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5781
$ws.Range("J131").Value = 3357.6
$ws.Range("L131").Value = 10072.8
$ws.Range("N131").Value = -20152.8

# Sheet GSM, row 107 (hunk idx 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 663.6923
$ws.Range("I107").Value = 529.875
$ws.Range("J107").Value = 877.8
$ws.Range("K107").Value = 529.875
$ws.Range("L107").Value = 877.8
$ws.Range("M107").Value = 1390.125
$ws.Range("N107").Value = -4717.8

# Sheet GSM, row 122 (hunk idx 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5153.579
$ws.Range("I122").Value = 5907.9287
$ws.Range("J122").Value = 3041.4
$ws.Range("K122").Value = 17723.7861
$ws.Range("L122").Value = 9124.200000000001
$ws.Range("M122").Value = -15273.7861
$ws.Range("N122").Value = -14024.2

# Sheet GSM, row 124 (hunk idx 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 60713.332
$ws.Range("J124").Value = 60713.332
$ws.Range("L124").Value = 60713.332
$ws.Range("N124").Value = -70533.33199999999

# Sheet LTW, row 122 (hunk idx 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2014.8422
$ws.Range("J122").Value = 2235.6667
$ws.Range("L122").Value = 6707.000100000001
$ws.Range("N122").Value = -11607.0001

# Sheet LTW, row 136 (hunk idx 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2446.4614
$ws.Range("I136").Value = 1178.2222
$ws.Range("J136").Value = 5300
$ws.Range("K136").Value = 3534.6666
$ws.Range("L136").Value = 15900
$ws.Range("M136").Value = -984.6665999999996
$ws.Range("N136").Value = -21000

# Sheet WVR, row 107 (hunk idx 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 631
$ws.Range("I107").Value = 401.83334
$ws.Range("J107").Value = 1456
$ws.Range("K107").Value = 1205.50002
$ws.Range("L107").Value = 4368
$ws.Range("M107").Value = 714.4999800000001
$ws.Range("N107").Value = -8208

# Sheet WVR, row 136 (hunk idx 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2701081.2
$ws.Range("I136").Value = 7412.4863
$ws.Range("J136").Value = 8930190
$ws.Range("K136").Value = 22237.4589
$ws.Range("L136").Value = 26790570
$ws.Range("M136").Value = -19687.4589
$ws.Range("N136").Value = -26795670
